# Scheduled market-data refresh: update price/profit figures across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6303.2856
$ws.Range("I62").Value = 5687.1665
$ws.Range("K62").Value = 5687.1665
$ws.Range("M62").Value = -5063.1665
$ws.Range("H65").Value = 6303.2856
$ws.Range("I65").Value = 5687.1665
$ws.Range("K65").Value = 28435.8325
$ws.Range("M65").Value = -25315.8325
$ws.Range("H88").Value = 20887738
$ws.Range("I88").Value = 111116250
$ws.Range("J88").Value = 65774
$ws.Range("K88").Value = 111116250
$ws.Range("L88").Value = 65774
$ws.Range("M88").Value = -111115844
$ws.Range("N88").Value = -66586
$ws.Range("H91").Value = 20887738
$ws.Range("I91").Value = 111116250
$ws.Range("J91").Value = 65774
$ws.Range("K91").Value = 111116250
$ws.Range("L91").Value = 65774
$ws.Range("M91").Value = -111114846
$ws.Range("N91").Value = -68582
$ws.Range("H113").Value = 25013992
$ws.Range("I113").Value = 2803.8
$ws.Range("J113").Value = 33351054
$ws.Range("K113").Value = 2803.8
$ws.Range("L113").Value = 33351054
$ws.Range("M113").Value = 450.1999999999998
$ws.Range("N113").Value = -33357562
$ws.Range("H132").Value = 1759.25
$ws.Range("I132").Value = 1759.25
$ws.Range("K132").Value = 5277.75
$ws.Range("M132").Value = -2747.75
$ws.Range("H138").Value = 4208.1465
$ws.Range("J138").Value = 4674.73
$ws.Range("L138").Value = 14024.19
$ws.Range("N138").Value = -24304.19

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3220.625
$ws.Range("I45").Value = 2748.8
$ws.Range("K45").Value = 2748.8
$ws.Range("M45").Value = -2371.8
$ws.Range("H97").Value = 3473526.5
$ws.Range("I97").Value = 917.9048
$ws.Range("J97").Value = 27781786
$ws.Range("K97").Value = 917.9048
$ws.Range("L97").Value = 27781786
$ws.Range("M97").Value = -421.9048
$ws.Range("N97").Value = -27782778
$ws.Range("H102").Value = 28577416
$ws.Range("I102").Value = 66671336
$ws.Range("J102").Value = 6974.75
$ws.Range("K102").Value = 66671336
$ws.Range("L102").Value = 6974.75
$ws.Range("M102").Value = -66669714
$ws.Range("N102").Value = -10218.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 2486.2856
$ws.Range("J14").Value = 3080.8
$ws.Range("L14").Value = 3080.8
$ws.Range("N14").Value = -3424.8
$ws.Range("H86").Value = 63746.35
$ws.Range("I86").Value = 95082.82000000001
$ws.Range("K86").Value = 95082.82000000001
$ws.Range("M86").Value = -93959.82000000001
$ws.Range("H89").Value = 63746.35
$ws.Range("I89").Value = 95082.82000000001
$ws.Range("K89").Value = 475414.1
$ws.Range("M89").Value = -469798.1
$ws.Range("H94").Value = 4873.8667
$ws.Range("I94").Value = 2153.4285
$ws.Range("J94").Value = 7254.25
$ws.Range("K94").Value = 2153.4285
$ws.Range("L94").Value = 7254.25
$ws.Range("M94").Value = -1702.4285
$ws.Range("N94").Value = -8156.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6758.7104
$ws.Range("I58").Value = 2229.0557
$ws.Range("J58").Value = 10835.4
$ws.Range("K58").Value = 2229.0557
$ws.Range("L58").Value = 10835.4
$ws.Range("M58").Value = -2026.0557
$ws.Range("N58").Value = -11241.4
$ws.Range("H105").Value = 5953797
$ws.Range("I105").Value = 7937286.5
$ws.Range("K105").Value = 7937286.5
$ws.Range("M105").Value = -7935539.5
$ws.Range("H107").Value = 1972.6538
$ws.Range("I107").Value = 1035.0769
$ws.Range("K107").Value = 1035.0769
$ws.Range("M107").Value = 884.9231
$ws.Range("H132").Value = 5728.5747
$ws.Range("I132").Value = 3210.0386
$ws.Range("K132").Value = 9630.1158
$ws.Range("M132").Value = -7100.1158
$ws.Range("H133").Value = 45200
$ws.Range("J133").Value = 50333.332
$ws.Range("L133").Value = 50333.332
$ws.Range("N133").Value = -55393.332
$ws.Range("H134").Value = 3796.0208
$ws.Range("I134").Value = 1114.5454
$ws.Range("K134").Value = 3343.6362
$ws.Range("M134").Value = -808.6361999999999
$ws.Range("H136").Value = 6758.7104
$ws.Range("I136").Value = 2229.0557
$ws.Range("J136").Value = 10835.4
$ws.Range("K136").Value = 6687.1671
$ws.Range("L136").Value = 32506.2
$ws.Range("M136").Value = -4137.1671
$ws.Range("N136").Value = -37606.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2153
$ws.Range("J5").Value = 3299.2144
$ws.Range("L5").Value = 9897.643199999999
$ws.Range("N5").Value = -10121.6432
$ws.Range("H14").Value = 18528586
$ws.Range("I14").Value = 18528586
$ws.Range("K14").Value = 55585758
$ws.Range("M14").Value = -55585585
$ws.Range("H26").Value = 612.1429000000001
$ws.Range("I26").Value = 408.33334
$ws.Range("J26").Value = 765
$ws.Range("K26").Value = 1225.00002
$ws.Range("L26").Value = 2295
$ws.Range("M26").Value = -937.0000199999999
$ws.Range("N26").Value = -2871
$ws.Range("H110").Value = 7291.5
$ws.Range("J110").Value = 8583.333000000001
$ws.Range("L110").Value = 25749.999
$ws.Range("N110").Value = -33929.999
$ws.Range("H131").Value = 38788.332
$ws.Range("I131").Value = 1786
$ws.Range("J131").Value = 54368.26
$ws.Range("K131").Value = 5358
$ws.Range("L131").Value = 163104.78
$ws.Range("M131").Value = -318
$ws.Range("N131").Value = -173184.78
$ws.Range("H135").Value = 2153
$ws.Range("J135").Value = 3299.2144
$ws.Range("L135").Value = 29692.9296
$ws.Range("N135").Value = -34762.9296

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 29999.5
$ws.Range("J101").Value = 29999.5
$ws.Range("L101").Value = 29999.5
$ws.Range("N101").Value = -36489.5
$ws.Range("H122").Value = 4027104.5
$ws.Range("I122").Value = 4529621.5
$ws.Range("K122").Value = 13588864.5
$ws.Range("M122").Value = -13586414.5
$ws.Range("H132").Value = 5412
$ws.Range("I132").Value = 2805.3125
$ws.Range("J132").Value = 12363.167
$ws.Range("K132").Value = 8415.9375
$ws.Range("L132").Value = 37089.501
$ws.Range("M132").Value = -5885.9375
$ws.Range("N132").Value = -42149.501

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4873.7812
$ws.Range("I40").Value = 4306.0386
$ws.Range("K40").Value = 4306.0386
$ws.Range("M40").Value = -4170.0386
$ws.Range("H132").Value = 11375074
$ws.Range("I132").Value = 19234890
$ws.Range("J132").Value = 22007
$ws.Range("K132").Value = 57704670
$ws.Range("L132").Value = 66021
$ws.Range("M132").Value = -57702140
$ws.Range("N132").Value = -71081
$ws.Range("H136").Value = 16577.162
$ws.Range("I136").Value = 5138.4
$ws.Range("K136").Value = 15415.2
$ws.Range("M136").Value = -12865.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 7247068
$ws.Range("I107").Value = 443.27585
$ws.Range("J107").Value = 19608958
$ws.Range("K107").Value = 1329.82755
$ws.Range("L107").Value = 58826874
$ws.Range("M107").Value = 590.17245
$ws.Range("N107").Value = -58830714
$ws.Range("H122").Value = 132615.64
$ws.Range("I122").Value = 201767.6
$ws.Range("J122").Value = 6884.8184
$ws.Range("K122").Value = 605302.8
$ws.Range("L122").Value = 20654.4552
$ws.Range("M122").Value = -602852.8
$ws.Range("N122").Value = -25554.4552
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 6167.3237
$ws.Range("I132").Value = 5653.88
$ws.Range("K132").Value = 16961.64
$ws.Range("M132").Value = -14431.64
